# Commerce_exterieur.xlsx correction
# - Replace the rich-text "Source : Direction Générale des Douanes/SYDONIA" note
#   with a plain-text " Direction Générale des Douanes/SYDONIA" string on the
#   sheets that show it (column E of "Imports_G_cont_t12" and "Alim val").
# - Refresh the remembered cell selection on every sheet, and make
#   "taux de change($)" the active sheet/tab (previously "BP " was active).

$wb = $excel.ActiveWorkbook

$newSource = ' Direction Générale des Douanes/SYDONIA'

# --- Imports_G_cont_t12 : column E (rows 2-62) holds the source note -------
$wsImports = $wb.Worksheets.Item('Imports_G_cont_t12')
for ($r = 2; $r -le 62; $r++) {
    $wsImports.Range("E$r").Value2 = $newSource
}

# --- Alim val : column E (rows 2-26) holds the source note ------------------
$wsAlim = $wb.Worksheets.Item('Alim val')
for ($r = 2; $r -le 26; $r++) {
    $wsAlim.Range("E$r").Value2 = $newSource
}

# --- Refresh remembered selections on every sheet ---------------------------

$wsTauxUSD = $wb.Worksheets.Item('taux de change($)')
$wsTauxUSD.Activate()
$wsTauxUSD.Range('I26').Select()

$wsTauxEUR = $wb.Worksheets.Item('taux de change(€)')
$wsTauxEUR.Activate()
$wsTauxEUR.Range('E2').Select()

$wsTauxFcfa = $wb.Worksheets.Item('taux de change(Fcfa)')
$wsTauxFcfa.Activate()
$wsTauxFcfa.Range('E2').Select()

$wsTauxYen = $wb.Worksheets.Item('taux de change(¥)')
$wsTauxYen.Activate()
$wsTauxYen.Range('E2').Select()

$wsImportBien = $wb.Worksheets.Item('Import bien')
$wsImportBien.Activate()
$wsImportBien.Range('E2').Select()

$wsImports.Activate()
$wsImports.Range('E2').Select()

$wsAlim.Activate()
$wsAlim.Range('M8').Select()

$wsExportTypr = $wb.Worksheets.Item('Exportation Mauritanienne ty.pr')
$wsExportTypr.Activate()
$wsExportTypr.Range('E2').Select()

$wsExportDes = $wb.Worksheets.Item('Exportation Mauritanien des')
$wsExportDes.Activate()
$wsExportDes.Range('E2').Select()

$wsExportPoisson = $wb.Worksheets.Item('Export poisson')
$wsExportPoisson.Activate()
$wsExportPoisson.Range('E2').Select()

$wsExportSmcp = $wb.Worksheets.Item('Export smcp pays')
$wsExportSmcp.Activate()
$wsExportSmcp.Range('E2').Select()

$wsBP = $wb.Worksheets.Item('BP ')
$wsBP.Activate()
$wsBP.Range('E2').Select()

# --- Make "taux de change($)" the active sheet/tab (was "BP ") -------------
$wsTauxUSD.Activate()
$wsTauxUSD.Range('I26').Select()
